# Added download order form feature:
# The "prices" sheet's online-storage pricing used a 2-tier (monthly/yearly
# per-GB) model. It is replaced with a 3-tier per-TB model:
#   online_storage_0_50_tb, online_storage_51_100_tb, online_storage_100_up_tb
# This pushes every row below the old online-storage rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prices")

# Insert a new row at position 12 - Excel shifts rows 12..40 down to 13..41
# and inherits the formatting of the row above (row 11, style s="27"),
# matching the target workbook's row 12 formatting.
$ws.Rows.Item(12).Insert()

# Row 10: first online-storage tier (0-50 TB) - plain (unstyled) numbers.
$ws.Cells.Item(10, 1).Value = "online_storage_0_50_tb"
$ws.Cells.Item(10, 2).Value = 0.048
$ws.Cells.Item(10, 3).Value = 0.033

# Row 11: second online-storage tier (51-100 TB). Overwrite the old label
# and values; this also clears the old "=0.04*12" formula in C11, replacing
# it with a plain literal value while keeping the existing number format.
$ws.Cells.Item(11, 1).Value = "online_storage_51_100_tb"
$ws.Cells.Item(11, 2).Value = 0.037
$ws.Cells.Item(11, 3).Value = 0.026

# Row 12 (newly inserted): third online-storage tier (100+ TB).
$ws.Cells.Item(12, 1).Value = "online_storage_100_up_tb"
$ws.Cells.Item(12, 2).Value = 0.028
$ws.Cells.Item(12, 3).Value = 0.019

# Update the view: scroll position / active selection moved.
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 5
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("D14").Select()
